$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-48.09999999999999, [double]"2.80738880071268e-06", -47.80829537787529, 90.82539594511528, 827.8255058766849, 0.7416168208955816, 0.4522805940438833, 36),
    @(-43.09999999999999, [double]"2.318331709998723e-06", -43.14635803198964, 91.0687371765809, 827.6612304733294, 0.8247158231892143, 0.7149604635334523, 50),
    @(-38.09999999999999, [double]"1.997579738166616e-06", -37.68518843851263, 90.78470942882738, 828.8377067354797, 0.6584141975665039, 0.7484436890083045, 66),
    @(-33.09999999999999, [double]"2.123849906631237e-06", -32.98388033521915, 90.90222861986571, 829.8173284104438, 0.9791104747640316, 0.6965302356205293, 93),
    @(-28.09999999999999, [double]"1.370316468192153e-06", -28.40995576350308, 90.65356287048206, 830.2837598993303, 0.793987156701334, 0.6679291429363076, 102),
    @(-23.09999999999999, [double]"1.336500239866959e-06", -22.86769729371633, 90.49122637440026, 830.6334606269025, 0.8056979721817132, 0.6697773104231175, 93),
    @(-18.09999999999999, [double]"6.339378657590692e-07", -18.20520397000868, 90.31453946716641, 830.7242974601622, 0.9503593598754119, 0.5999388550668032, 109),
    @(-13.09999999999999, [double]"-3.923298197014747e-07", -13.53538203551962, 90.40648246556162, 831.0423878309792, 1.291263089249161, 0.5191800489044687, 118),
    @(-8.099999999999994, [double]"-2.201736949984444e-07", -8.412366102243347, 90.13248741757634, 830.8335811319689, 1.454014770682649, 0.6738613806618279, 130),
    @(-3.099999999999994, [double]"-1.780561801267056e-06", -4.162115875665234, 186.8533940033314, 830.5485616267192, 1.92309375837692, 0.6170306428121161, 110),
    @(1.900000000000006, [double]"-1.45642546286813e-06", 1.5, 89.70839706418597, 771.5078743507485, 2.317933490908868, 0.5909648256944725, 147),
    @(6.900000000000006, [double]"-2.597801391765845e-06", 6.883590571785827, 90.01028177563565, 771.6986054665356, 0.9630759392654894, 0.6927977533583594, 108),
    @(11.90000000000001, [double]"-3.972583130352006e-06", 12.29422759822129, 90.13215442508441, 771.9591986155116, 0.8578763283534891, 0.7300137378305677, 113),
    @(16.90000000000001, [double]"-4.150207436159695e-06", 17.00708776609423, 90.11881280976174, 830.640710974608, 0.9094950290507534, 0.6966699915145911, 121),
    @(21.90000000000001, [double]"-4.741710607907295e-06", 22.14348712706708, 90.2179478546225, 830.2221608451788, 0.8859619568781814, 0.744676304663677, 129),
    @(26.90000000000001, [double]"-5.35908295853718e-06", 27.84518183438559, 90.18888841688099, 829.7983117699556, 0.8296708479267555, 0.7849704397612036, 110),
    @(31.90000000000001, [double]"-5.749188643935305e-06", 32.63277172314476, 90.15951874917874, 772.4037783531151, 0.8516371813098846, 0.8454017226245476, 105),
    @(36.90000000000001, [double]"-6.353637337595844e-06", 37.69291482618024, 90.29657357722289, 773.0103223849046, 1.018444302602618, 0.8347644503624119, 99),
    @(41.90000000000001, [double]"-5.173322218226858e-06", 42.43896635786598, 89.88634796144834, 772.4977024774959, 1.054231890990177, 0.7408134199985423, 59),
    @(46.90000000000001, [double]"-4.696845708077501e-06", 47.73905463312075, 90.4223439354706, 772.6022991595236, 0.8897300773309664, 0.6814878313815504, 42),
    @(51.90000000000001, [double]"-5.169729113274791e-06", 53.27214888613018, 90.54075027739067, 772.2290571380361, 1.072497289392378, 0.8369189191930244, 37)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
